# Auto-generated edit script applying diff changes to Jogos_do_Dia_Betfair_Back_Lay_2025-12-18.xlsx
# Updates numeric odds values in rows 2-19 of Sheet1, per the committed diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("H2").Value = 5
$ws.Range("J2").Value = 4
$ws.Range("P2").Value = 2.2
$ws.Range("Q2").Value = 1.63
$ws.Range("Z2").Value = 55
$ws.Range("AI2").Value = 75
# Row 3
$ws.Range("H3").Value = 3.15
$ws.Range("X3").Value = 18.5
$ws.Range("Y3").Value = 15.5
$ws.Range("AA3").Value = 60
$ws.Range("AC3").Value = 9.199999999999999
$ws.Range("AD3").Value = 14.5
$ws.Range("AE3").Value = 36
$ws.Range("AG3").Value = 12
$ws.Range("AH3").Value = 17
$ws.Range("AI3").Value = 44
$ws.Range("AM3").Value = 80
$ws.Range("AN3").Value = 17
# Row 4
$ws.Range("F4").Value = 1.26
$ws.Range("G4").Value = 1.28
$ws.Range("H4").Value = 12.5
$ws.Range("I4").Value = 17
$ws.Range("K4").Value = 7.8
$ws.Range("N4").Value = 4.8
$ws.Range("O4").Value = 1.16
$ws.Range("P4").Value = 2.76
$ws.Range("Q4").Value = 1.56
$ws.Range("R4").Value = 1.68
$ws.Range("T4").Value = 2.16
$ws.Range("U4").Value = 1.68
$ws.Range("X4").Value = 34
$ws.Range("Y4").Value = 55
$ws.Range("Z4").Value = 170
$ws.Range("AB4").Value = 13.5
$ws.Range("AC4").Value = 19
$ws.Range("AD4").Value = 65
$ws.Range("AE4").Value = 310
$ws.Range("AF4").Value = 8.6
$ws.Range("AG4").Value = 14.5
$ws.Range("AH4").Value = 40
$ws.Range("AI4").Value = 210
$ws.Range("AJ4").Value = 10
$ws.Range("AK4").Value = 15.5
$ws.Range("AM4").Value = 210
$ws.Range("AN4").Value = 4.3
# Row 5
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 5.8
$ws.Range("Q5").Value = 2.06
# Row 6
$ws.Range("N6").Value = 3.4
$ws.Range("T6").Value = 1.79
$ws.Range("U6").Value = 2.02
$ws.Range("X6").Value = 17.5
$ws.Range("Y6").Value = 18.5
$ws.Range("AC6").Value = 9.800000000000001
# Row 7
$ws.Range("H7").Value = 5.1
$ws.Range("K7").Value = 7
$ws.Range("P7").Value = 2.24
$ws.Range("Q7").Value = 1.53
$ws.Range("T7").Value = 2.08
$ws.Range("U7").Value = 1.76
$ws.Range("Y7").Value = 44
$ws.Range("Z7").Value = 140
$ws.Range("AC7").Value = 980
$ws.Range("AF7").Value = 8.4
$ws.Range("AG7").Value = 11.5
$ws.Range("AJ7").Value = 11
$ws.Range("AK7").Value = 1000
$ws.Range("AN7").Value = 5.3
# Row 8
$ws.Range("I8").Value = 4.2
$ws.Range("J8").Value = 3.7
$ws.Range("K8").Value = 3.85
$ws.Range("P8").Value = 1.69
$ws.Range("Q8").Value = 1.99
$ws.Range("T8").Value = 1.96
$ws.Range("U8").Value = 1.86
$ws.Range("AA8").Value = 140
$ws.Range("AM8").Value = 180
# Row 9
$ws.Range("F9").Value = 2.54
$ws.Range("T9").Value = 1.82
$ws.Range("U9").Value = 1.99
$ws.Range("AB9").Value = 10.5
$ws.Range("AD9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AM9").Value = 130
# Row 10
$ws.Range("H10").Value = 2.74
$ws.Range("I10").Value = 2.84
$ws.Range("K10").Value = 2.88
# Row 11
$ws.Range("F11").Value = 1.24
$ws.Range("H11").Value = 1.5
$ws.Range("K11").Value = 8.6
$ws.Range("P11").Value = 2.76
$ws.Range("T11").Value = 1.87
$ws.Range("U11").Value = 1.91
$ws.Range("AB11").Value = 1000
$ws.Range("AC11").Value = 22
$ws.Range("AH11").Value = 38
$ws.Range("AI11").Value = 140
$ws.Range("AJ11").Value = 11
$ws.Range("AN11").Value = 3.65
# Row 12
$ws.Range("F12").Value = 1.29
$ws.Range("G12").Value = 1.31
$ws.Range("H12").Value = 8.4
$ws.Range("J12").Value = 5.8
$ws.Range("K12").Value = 7
$ws.Range("P12").Value = 2.36
$ws.Range("T12").Value = 2.04
$ws.Range("AB12").Value = 10
$ws.Range("AF12").Value = 9
$ws.Range("AK12").Value = 1000
$ws.Range("AM12").Value = 200
$ws.Range("AN12").Value = 5.1
# Row 13
$ws.Range("F13").Value = 1.45
$ws.Range("H13").Value = 6.8
$ws.Range("K13").Value = 5.2
$ws.Range("P13").Value = 1.89
$ws.Range("Q13").Value = 1.77
$ws.Range("T13").Value = 2.14
$ws.Range("U13").Value = 1.7
$ws.Range("X13").Value = 16.5
$ws.Range("AA13").Value = 490
$ws.Range("AB13").Value = 8.6
$ws.Range("AC13").Value = 11.5
$ws.Range("AD13").Value = 46
$ws.Range("AE13").Value = 230
$ws.Range("AF13").Value = 9.199999999999999
$ws.Range("AG13").Value = 12.5
$ws.Range("AH13").Value = 32
$ws.Range("AI13").Value = 190
$ws.Range("AJ13").Value = 14
$ws.Range("AK13").Value = 20
$ws.Range("AL13").Value = 50
$ws.Range("AM13").Value = 250
$ws.Range("AN13").Value = 8.4
# Row 14
$ws.Range("H14").Value = 1.09
$ws.Range("I14").Value = 34
$ws.Range("J14").Value = 8
$ws.Range("P14").Value = 2.58
$ws.Range("Q14").Value = 1.48
$ws.Range("T14").Value = 2.4
$ws.Range("AH14").Value = 50
$ws.Range("AJ14").Value = 8.6
$ws.Range("AL14").Value = 55
# Row 15
$ws.Range("F15").Value = 2.16
$ws.Range("I15").Value = 3.7
$ws.Range("P15").Value = 1.79
$ws.Range("Q15").Value = 1.74
$ws.Range("U15").Value = 2.1
$ws.Range("AE15").Value = 1000
$ws.Range("AF15").Value = 15
$ws.Range("AH15").Value = 21
$ws.Range("AL15").Value = 1000
$ws.Range("AO15").Value = 1000
# Row 16
$ws.Range("P16").Value = 2.94
$ws.Range("Q16").Value = 1.41
$ws.Range("Y16").Value = 55
$ws.Range("AB16").Value = 13
$ws.Range("AF16").Value = 10.5
$ws.Range("AN16").Value = 3.75
# Row 17
$ws.Range("F17").Value = 3.4
$ws.Range("AC17").Value = 10
# Row 18
$ws.Range("F18").Value = 1.9
$ws.Range("G18").Value = 2.02
$ws.Range("J18").Value = 3.45
$ws.Range("X18").Value = 1000
$ws.Range("AF18").Value = 1000
# Row 19
$ws.Range("F19").Value = 1.7
$ws.Range("H19").Value = 4.5
$ws.Range("K19").Value = 4.6
$ws.Range("Q19").Value = 1.63
$ws.Range("T19").Value = 1.77
$ws.Range("U19").Value = 2.04
$ws.Range("X19").Value = 22
$ws.Range("AA19").Value = 170
$ws.Range("AB19").Value = 11.5
$ws.Range("AJ19").Value = 21
$ws.Range("AK19").Value = 21
$ws.Range("AN19").Value = 11
